$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Formula = "=AVERAGE(A2:A4)"
$ws.Range("B5:K5").Formula = "=AVERAGE(B2:B4)"

$ws.Range("A5:K5").Font.Bold = $true

$ws.Range("O5").Select()
